$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 15: scale the "curve 1" reference row (row 4) currents/voltages by
# 5x (currents, B:G) and 12x (voltages, H:L) to derive a new parameter row.
$ws.Range("B15").Formula = "=B4*5"
$ws.Range("C15").Formula = "=C4*5"
$ws.Range("D15").Formula = "=D4*5"
$ws.Range("E15").Formula = "=E4*5"
$ws.Range("F15").Formula = "=F4*5"
$ws.Range("G15").Formula = "=G4*5"
$ws.Range("H15").Formula = "=H4*12"
$ws.Range("I15:L15").Formula = "=I4*12"

$ws.Range("J15").Select()
